$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the three header cells in row 5 and a few others that were using
# the (border+fill) style; the "fill" was never actually visible (patternType
# none) so the workbook is being cleaned up to just use the plain bordered
# style that every other cell already uses. Do this by copying the format
# from a cell that already has the correct (border-only) style.
$borderOnly = $ws.Range("B5")
$borderOnly.Copy()
foreach ($addr in @("H5", "I5", "J5", "J7", "J8", "B11", "B12")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- New "Поля юзера" (user fields) table, mirroring the existing
# "Поля комнаты" table layout: a header row followed by one data row.
$ws.Range("B28").Value = "Поля юзера"
$ws.Range("C28").Value = "telegram_ID"
$ws.Range("D28").Value = "Имя"
$ws.Range("E28").Value = "группа"
$ws.Range("F28").Value = "роль в группе"

$ws.Range("B29").Value = "Значения по умолчанию из класса"
$ws.Range("C29").Value = "да"
$ws.Range("D29").Value = "да"
$ws.Range("E29").Value = "да"
$ws.Range("F29").Value = "да"

# Give the new rows the same bordered look as the rest of the sheet.
$borderOnly.Copy()
$ws.Range("B28:F29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the saved selection to match where the author left off editing.
$ws.Range("I34").Select()
